$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ResetPassword")

# 1. The generic "Error" message scenarios now describe the actual error text
#    returned by the UI when a reset-password token is invalid/expired.
$ws.Range("E3").Value = "Invalid or expired reset token"
$ws.Range("E7").Value = "Invalid or expired reset token"

# 2. Add the missing "happy path" scenario: submitting the reset-password
#    form with all valid inputs results in a success message being shown.
$ws.Range("A8").Value = "Submit with all Valid Inputs"
$ws.Range("C8").Value = "password"
$ws.Range("D8").Value = "password"
$ws.Range("E8").Value = "Your password has been reset successfully"

# Match the formatting used by the other scenario/message cells in the sheet
# (copy the look of the equivalent column from a neighbouring row).
$ws.Range("A6").Copy()
$ws.Range("A8").PasteSpecial(-4122)

$ws.Range("C6").Copy()
$ws.Range("C8").PasteSpecial(-4122)

$ws.Range("D6").Copy()
$ws.Range("D8").PasteSpecial(-4122)

$ws.Range("E4").Copy()
$ws.Range("E8").PasteSpecial(-4122)

$ws.Activate()
$ws.Range("E6").Select()
